$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.565.89"
$ws.Range("E2").Value = "  +2.99%  "
$ws.Range("D3").Value = "1.695.43"
$ws.Range("E3").Value = "  +2.61%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.83%  "
$ws.Range("D5").Value = "'313.34"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.70%  "
$ws.Range("D6").Value = "'1.000"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Value = "'0.3952"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.79%  "
$ws.Range("B8").Value = "OKB"
$ws.Range("C8").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D8").Value = "'58.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +21.65%  "
$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D9").Value = "'0.4048"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.50%  "
$ws.Range("D10").Value = "'1.531"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +11.27%  "
$ws.Range("D11").Value = "'1.002"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.79%  "
$ws.Range("D12").Value = "'0.08778"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.59%  "
$ws.Range("D13").Value = "'7.293"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +14.14%  "
$ws.Range("E14").Value = "  +3.36%  "
$ws.Range("D15").Value = "'0.00001319"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.74%  "
$ws.Range("D16").Value = "'7.655"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +7.98%  "
$ws.Range("D17").Value = "1.697.85"
$ws.Range("E17").Value = "  +2.32%  "
$ws.Range("D18").Value = "'100.43"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.63%  "
$ws.Range("D19").Value = "'0.07051"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.47%  "
$ws.Range("E20").Value = "  +4.05%  "
$ws.Range("D21").Value = "'6.739"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.06%  "
$ws.Range("D22").Value = "'1.000"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.71%  "
$ws.Range("D23").Value = "'14.18"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.85%  "
$ws.Range("D24").Value = "24.557.42"
$ws.Range("E24").Value = "  +2.93%  "
$ws.Range("D25").Value = "'3.009"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +14.25%  "
$ws.Range("D26").Value = "'2.309"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.16%  "
$ws.Range("D27").Value = "'22.50"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.27%  "
$ws.Range("D28").Value = "'159.78"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.67%  "
$ws.Range("D29").Value = "'5.178"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.13%  "
$ws.Range("D30").Value = "'133.67"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.52%  "
$ws.Range("D31").Value = "'7.682"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +40.70%  "
$ws.Range("D32").Value = "1.883.67"
$ws.Range("E32").Value = "  +2.24%  "
$ws.Range("D33").Value = "'1.091"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.01%  "
$ws.Range("D34").Value = "'7.321"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +22.07%  "
$ws.Range("D35").Value = "'0.08534"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.21%  "
$ws.Range("D36").Value = "'1.960"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +10.87%  "
$ws.Range("D37").Value = "'11.00"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +7.28%  "
$ws.Range("D38").Value = "'0.2720"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.98%  "
$ws.Range("D39").Value = "'14.72"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.70%  "
$ws.Range("D40").Value = "'0.02783"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +12.87%  "
$ws.Range("D41").Value = "'0.09065"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.70%  "
$ws.Range("D42").Value = "'1.472"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.47%  "
$ws.Range("D43").Value = "'0.7638"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.04%  "
$ws.Range("D44").Value = "'0.7176"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.73%  "
$ws.Range("D45").Value = "'15.33"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.13%  "
$ws.Range("D46").Value = "'2.462"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.66%  "
$ws.Range("D47").Value = "'4.178"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.06%  "
$ws.Range("D48").Value = "'0.9991"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.79%  "
$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").Value = "'140.45"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.94%  "
$ws.Range("B50").Value = "Flow"
$ws.Range("C50").Value = "https://coinranking.com/coin/QQ0NCmjVq+flow-flow"
$ws.Range("D50").Value = "'1.321"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +18.98%  "
$ws.Range("E51").Value = "  +2.04%  "
